$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Range("D7").Value = "INT2204 1"
$ws.Range("E7").Value = "Hương đối tượng"
$ws.Range("G7").Value = 1
$ws.Range("H7").Value = "Tô Văn Khánh"

$ws.Range("H8").Select()
